$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22 (Fri, Oct 13): rename the midterm review lab entry and add its
# reading link, consistent with the new "Lab 6 / Lab 7" naming below.
$ws.Range("C22").Value = "Lab : Review for Midterm I"
$ws.Range("F22").Value = "lab-midterm-I-review.html"

# Row 25 (Fri, Oct 20): this used to be the placeholder "Lab 7: " entry
# with the HW7 link and "Lab 7: Variable Selection" in the Old column.
# It becomes "Lab 6: Dimension Reduction via Regularization", the HW7
# link moves down to row 26, and the Old-column note is renamed to
# "Lab 7: Dimension Reduction via Regularization".
$ws.Range("C25").Value = "Lab 6: Dimension Reduction via Regularization"
$ws.Range("G25").Value = ""
$ws.Range("H25").Value = "Lab 7: Dimension Reduction via Regularization"

# Row 26 (WEEK 9 / Tue, Oct 24): receives the HW7 link that moved off row 25.
$ws.Range("G26").Value = "hw-07"

# Update the saved selection to match the authored state.
$null = $ws.Range("C22").Select()
